$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# 1. Update the three (Kyrgyz/Russian/English) title cells in row 1:
#    "4.c.1.1 ..." -> "4.c.1 ..." (indicator code correction)
# ------------------------------------------------------------------
$ws.Range("A1").Value = "4.с.1 Билим берүү мекемелерде диплом берилгем мугалимдердин үлүшү"
$ws.Range("B1").Value = "4.c.1 Доля дипломированных учителей в образовательных учереждениях"
$ws.Range("C1").Value = "4.c.1 Proportion of certified teachers in educational institutions"

# ------------------------------------------------------------------
# 2. Remember the formatting of the cells that currently sit in the
#    last data column (M) and at L8 -- we need these *before* the
#    column layout changes so we can re-apply them afterwards.
# ------------------------------------------------------------------
$ws.Range("M4").Copy() | Out-Null
$ws.Range("ZZ1").PasteSpecial(-4122) | Out-Null   # xlPasteFormats -> stash year-header format

$ws.Range("M5").Copy() | Out-Null
$ws.Range("ZZ2").PasteSpecial(-4122) | Out-Null   # stash plain data-row format (rows 5 & 6)

$ws.Range("M7").Copy() | Out-Null
$ws.Range("ZZ3").PasteSpecial(-4122) | Out-Null   # stash bordered data-row format (row 7)

$ws.Range("L8").Copy() | Out-Null
$ws.Range("ZZ4").PasteSpecial(-4122) | Out-Null   # stash the little bold marker cell's format

# ------------------------------------------------------------------
# 3. Insert columns for the new years: 2011 & 2012 (right after 2010)
#    and 2022 (right after 2021).
# ------------------------------------------------------------------
$ws.Range("E:F").Insert() | Out-Null
$ws.Range("P:P").Insert() | Out-Null

# Clear whatever the insert/shift operation left behind in row 8 past
# column L (the little marker cell must only exist at L8).
$ws.Range("M8:P8").Clear() | Out-Null

# ------------------------------------------------------------------
# 4. Fill in the year header row and the three data rows with the
#    complete 2010-2022 series.
# ------------------------------------------------------------------
$ws.Range("D4").Value = 2010
$ws.Range("E4").Value = 2011
$ws.Range("F4").Value = 2012
$ws.Range("G4").Value = 2013
$ws.Range("H4").Value = 2014
$ws.Range("I4").Value = 2015
$ws.Range("J4").Value = 2016
$ws.Range("K4").Value = 2017
$ws.Range("L4").Value = 2018
$ws.Range("M4").Value = 2019
$ws.Range("N4").Value = 2020
$ws.Range("O4").Value = 2021
$ws.Range("P4").Value = 2022

$ws.Range("D5").Value = 87.9
$ws.Range("E5").Value = 89.6
$ws.Range("F5").Value = 87.5
$ws.Range("G5").Value = 88.8
$ws.Range("H5").Value = 89.8
$ws.Range("I5").Value = 94.7
$ws.Range("J5").Value = 91.6
$ws.Range("K5").Value = 93.4
$ws.Range("L5").Value = 93.5
$ws.Range("M5").Value = 93.6
$ws.Range("N5").Value = 94.5
$ws.Range("O5").Value = 93.5
$ws.Range("P5").Value = 94.2

$ws.Range("D6").Value = 93.6
$ws.Range("E6").Value = 93.3
$ws.Range("F6").Value = 93.9
$ws.Range("G6").Value = 94.3
$ws.Range("H6").Value = 94.4
$ws.Range("I6").Value = 95
$ws.Range("J6").Value = 95.4
$ws.Range("K6").Value = 96
$ws.Range("L6").Value = 96.4
$ws.Range("M6").Value = 96.3
$ws.Range("N6").Value = 96.7
$ws.Range("O6").Value = 96.6
$ws.Range("P6").Value = 96

$ws.Range("D7").Value = 92.9
$ws.Range("E7").Value = 92.8
$ws.Range("F7").Value = 94.1
$ws.Range("G7").Value = 94.8
$ws.Range("H7").Value = 95.3
$ws.Range("I7").Value = 95.9
$ws.Range("J7").Value = 96.9
$ws.Range("K7").Value = 97.9
$ws.Range("L7").Value = 98
$ws.Range("M7").Value = 98
$ws.Range("N7").Value = 98.2
$ws.Range("O7").Value = 98.1
$ws.Range("P7").Value = 97.5

# ------------------------------------------------------------------
# 5. Re-apply the stashed formats across the full, now-wider ranges.
# ------------------------------------------------------------------
$ws.Range("ZZ1").Copy() | Out-Null
$ws.Range("D4:P4").PasteSpecial(-4122) | Out-Null

$ws.Range("ZZ2").Copy() | Out-Null
$ws.Range("D5:P5").PasteSpecial(-4122) | Out-Null
$ws.Range("D6:P6").PasteSpecial(-4122) | Out-Null

$ws.Range("ZZ3").Copy() | Out-Null
$ws.Range("D7:P7").PasteSpecial(-4122) | Out-Null

$ws.Range("ZZ4").Copy() | Out-Null
$ws.Range("L8").PasteSpecial(-4122) | Out-Null

# P6 (2022 value for the "basic school" row) is a whole number but is
# displayed with one decimal place to match the rest of the series.
$ws.Range("P6").NumberFormat = "0.0"

# ------------------------------------------------------------------
# 6. Clean up the scratch area used to stash formats.
# ------------------------------------------------------------------
$ws.Range("ZZ1:ZZ4").Clear() | Out-Null
$excel.CutCopyMode = 0

# ------------------------------------------------------------------
# 7. Cosmetic: move the active selection the way the source file has it.
# ------------------------------------------------------------------
$ws.Range("Q4").Select() | Out-Null

Write-Output "done"
